# Regenerate sval data to filter save games.
# Updates the B:E (and derived G=sum) columns for rows 2-12 on the active sheet
# with the new, recomputed statistics, leaving column F (Win) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 13.86384647080068; G = 19.48425592650926 }
    3  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    4  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    5  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3993.344853322108;  E = 13.86384647080068; G = 4012.10801473063 }
    6  = @{ B = 0.6545652718822623; C = 0.3048912486333797; D = 0.1496068669990043; E = 0.5333859586016987; G = 1.642449346116345 }
    7  = @{ B = 0.2881169905109251; C = 0.3048912486333797; D = 3.223369029078222;  E = 0.5333859586016987; G = 4.349763226824225 }
    8  = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    9  = @{ B = 0.6545652718822623; C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 6.038307959104277 }
    10 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    11 = @{ B = 0.01253208636536152; C = 1.626987699542094; D = 3.223369029078222; E = 13.86384647080068; G = 18.72673528578636 }
    12 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
